{"js": "const replacements = [\n  [\"2025-11-24 Monday\", \"2025-11-25 Tuesday\"],\n  [\"14\u00d711=154\", \"60\u00d748=2880\"],\n  [\"95\u00d798=9310\", \"35\u00d752=1820\"],\n  [\"88\u00d744=3872\", \"35\u00d724=840\"],\n  [\"96\u00d767=6432\", \"58\u00d793=5394\"],\n  [\"90\u00d738=3420\", \"60\u00d745=2700\"],\n  [\"95\u00d796=9120\", \"48\u00d716=768\"],\n  [\"17\u00d721=357\", \"97\u00d778=7566\"],\n  [\"53\u00d738=2014\", \"53\u00d762=3286\"],\n  [\"59\u00d772=4248\", \"55\u00d776=4180\"],\n  [\"99\u00d772=7128\", \"96\u00d761=5856\"],\n  [\"27\u00d730=810\", \"79\u00d755=4345\"],\n  [\"42\u00d772=3024\", \"89\u00d712=1068\"],\n  [\"84\u00d712=1008\", \"40\u00d736=1440\"],\n  [\"79\u00d724=1896\", \"36\u00d736=1296\"],\n  [\"27\u00d738=1026\", \"75\u00d721=1575\"],\n  [\"92\u00d787=8004\", \"36\u00d735=1260\"],\n  [\"54\u00d716=864\", \"72\u00d717=1224\"],\n  [\"48\u00d759=2832\", \"75\u00d791=6825\"],\n  [\"41\u00d793=3813\", \"23\u00d788=2024\"],\n  [\"61\u00d715=915\", \"55\u00d763=3465\"],\n  [\"89\u00d787=7743\", \"11\u00d713=143\"],\n  [\"67\u00d760=4020\", \"17\u00d754=918\"],\n  [\"14\u00d780=1120\", \"83\u00d752=4316\"],\n  [\"12\u00d790=1080\", \"24\u00d790=2160\"],\n  [\"40\u00d733=1320\", \"68\u00d790=6120\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute('2025-11-24 Monday', $false, $false, $false, $false, $false, $true, 1, $false, '2025-11-25 Tuesday', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute('14\u00d711=154', $false, $false, $false, $false, $false, $true, 1, $false, '60\u00d748=2880', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute('95\u00d798=9310', $false, $false, $false, $false, $false, $true, 1, $false, '35\u00d752=1820', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute('88\u00d744=3872', $false, $false, $false, $false, $false, $true, 1, $false, '35\u00d724=840', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute('96\u00d767=6432', $false, $false, $false, $false, $false, $true, 1, $false, '58\u00d793=5394', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute('90\u00d738=3420', $false, $false, $false, $false, $false, $true, 1, $false, '60\u00d745=2700', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute('95\u00d796=9120', $false, $false, $false, $false, $false, $true, 1, $false, '48\u00d716=768', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute('17\u00d721=357', $false, $false, $false, $false, $false, $true, 1, $false, '97\u00d778=7566', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute('53\u00d738=2014', $false, $false, $false, $false, $false, $true, 1, $false, '53\u00d762=3286', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute('59\u00d772=4248', $false, $false, $false, $false, $false, $true, 1, $false, '55\u00d776=4180', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute('99\u00d772=7128', $false, $false, $false, $false, $false, $true, 1, $false, '96\u00d761=5856', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute('27\u00d730=810', $false, $false, $false, $false, $false, $true, 1, $false, '79\u00d755=4345', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute('42\u00d772=3024', $false, $false, $false, $false, $false, $true, 1, $false, '89\u00d712=1068', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute('84\u00d712=1008', $false, $false, $false, $false, $false, $true, 1, $false, '40\u00d736=1440', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute('79\u00d724=1896', $false, $false, $false, $false, $false, $true, 1, $false, '36\u00d736=1296', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute('27\u00d738=1026', $false, $false, $false, $false, $false, $true, 1, $false, '75\u00d721=1575', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute('92\u00d787=8004', $false, $false, $false, $false, $false, $true, 1, $false, '36\u00d735=1260', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute('54\u00d716=864', $false, $false, $false, $false, $false, $true, 1, $false, '72\u00d717=1224', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute('48\u00d759=2832', $false, $false, $false, $false, $false, $true, 1, $false, '75\u00d791=6825', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute('41\u00d793=3813', $false, $false, $false, $false, $false, $true, 1, $false, '23\u00d788=2024', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute('61\u00d715=915', $false, $false, $false, $false, $false, $true, 1, $false, '55\u00d763=3465', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute('89\u00d787=7743', $false, $false, $false, $false, $false, $true, 1, $false, '11\u00d713=143', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute('67\u00d760=4020', $false, $false, $false, $false, $false, $true, 1, $false, '17\u00d754=918', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute('14\u00d780=1120', $false, $false, $false, $false, $false, $true, 1, $false, '83\u00d752=4316', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute('12\u00d790=1080', $false, $false, $false, $false, $false, $true, 1, $false, '24\u00d790=2160', 2) | Out-Null\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute('40\u00d733=1320', $false, $false, $false, $false, $false, $true, 1, $false, '68\u00d790=6120', 2) | Out-Null\n"}
